$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24: pick up what used to be row 25's survey data ---
$ws.Range("A24").Value = 111573746
$ws.Range("Q24").Value = 562601.7570288588
$ws.Range("R24").Value = 6954814.918206804
$ws.Range("Z24").Value = "15:12"
$ws.Range("AB24").Value = "15:12"

# --- Row 25: pick up what used to be row 26's survey data ---
$ws.Range("A25").Value = 111575785
$ws.Range("B25").Value = 89845
$ws.Range("E25").Value = 1209
$ws.Range("F25").Value = "Rynkskinn"
$ws.Range("G25").Value = "Phlebia centrifuga"
$ws.Range("H25").Value = "P.Karst."
$ws.Range("Q25").Value = 562859.2727272335
$ws.Range("R25").Value = 6954660.134623887
$ws.Range("Z25").Value = "16:39"
$ws.Range("AB25").Value = "16:39"

# --- Row 26: pick up what used to be row 24's survey data ---
$ws.Range("A26").Value = 111578197
$ws.Range("B26").Value = 96348
$ws.Range("E26").Value = 220787
$ws.Range("F26").Value = "Knärot"
$ws.Range("G26").Value = "Goodyera repens"
$ws.Range("H26").Value = "(L.) R. Br."
$ws.Range("Q26").Value = 563026.0554397166
$ws.Range("R26").Value = 6954541.256262898
$ws.Range("Z26").Value = "00:00"
$ws.Range("AB26").Value = "00:00"

# --- Row 32 and Row 33 swap their survey data ---
$ws.Range("A32").Value = 111576450
$ws.Range("Q32").Value = 562979.5212303887
$ws.Range("R32").Value = 6954739.97881452
$ws.Range("Z32").Value = "17:10"
$ws.Range("AB32").Value = "17:10"
$ws.Range("AC32").Value = "Rikligt"

$ws.Range("A33").Value = 111576771
$ws.Range("Q33").Value = 562807.4867926922
$ws.Range("R33").Value = 6954821.585021482
$ws.Range("Z33").Value = "17:24"
$ws.Range("AB33").Value = "17:24"
$ws.Range("AC33").ClearContents()
